$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.194.59'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.181.42'
$ws.Range('E3').Value = '  -2.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.43'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -7.27%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.572'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.19'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -13.08%  '
$ws.Range('E12').Value = '  -4.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.104'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.56%  '
$ws.Range('D15').Value = '2.508.57'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.40'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.855'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '2.209.38'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').Value = '41.195.82'
$ws.Range('E19').Value = '  -1.51%  '
$ws.Range('D20').Value = '0.0₃0947'
$ws.Range('E20').Value = '  -2.75%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.68'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('E24').Value = '  -6.06%  '
$ws.Range('E25').Value = '  -5.68%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.91%  '
$ws.Range('E28').Value = '  -4.51%  '
$ws.Range('E29').Value = '  -4.11%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.11'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '168.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.123'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0735'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.122'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.99%  '
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.43'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0308'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.54%  '
$ws.Range('E41').Value = '  -3.44%  '
$ws.Range('E42').Value = '  +7.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.31%  '
$ws.Range('E45').Value = '  -5.89%  '
$ws.Range('E46').Value = '  -2.58%  '
$ws.Range('E47').Value = '  -8.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1000'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('E51').Value = '  -3.24%  '
